$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "18.8.2019"
$ws.Range("B21").Value = 12
$ws.Range("C21").Value = "Karmeasti frontin refaktorointia: komponenttihakemistojen rakenteen pohtimista ja järjestelyä, TypeScriptin oikean käytön selvittelyä ja tyyppimäärittelyjen tekoa uudelleen, Redux storen iso refaktorointi. Git katastrofin korjausta 1 h."

$ws.Rows.Item(21).RowHeight = 77.3
